$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Collapse proofing-split runs back into single runs (no visible text
#    change, just removes the spell-check run-splitting / proofErr markers).
#    A self Find & Replace on the whole-paragraph text forces Word to
#    re-emit the paragraph as a single run.
# ---------------------------------------------------------------------------

$mergeTexts = @(
    "Willen hun ticket per e-mail in PDF-formaat ontvangen.",
    "Beveiligde online betalingsopties (iDEAL, creditcard, PayPal).",
    "Responsieve website die werkt op mobiele apparaten en desktops.",
    "Cloudgebaseerde opslag van ticketgegevens voor toegankelijkheid."
)

foreach ($t in $mergeTexts) {
    $d.Content.Find.Execute($t, $true, $true, $false, $false, $false, $true, 1, $false, $t, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Re-order the bullet list under "3.1 Basisfunctionaliteiten":
#    remove "Tickets worden als PDF per e-mail verstuurd." and shift the
#    two following bullets up, then append "Meertalige ondersteuning
#    (Nederlands en Engels)." followed by a new empty bullet.
# ---------------------------------------------------------------------------

$d.Paragraphs.Item(20).Range.Text = "Een inlogsysteem voor medewerkers."
$d.Paragraphs.Item(21).Range.Text = "Mogelijkheid voor medewerkers om tickets te scannen en valideren."
$d.Paragraphs.Item(22).Range.Text = "Meertalige ondersteuning (Nederlands en Engels)."

$tailRange = $d.Paragraphs.Item(22).Range
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 3) Remove the now-duplicate "Meertalige ondersteuning (Nederlands en
#    Engels)." bullet from "3.2 Advanced Functionaliteiten" (it has moved
#    to 3.1 above).
# ---------------------------------------------------------------------------

$d.Paragraphs.Item(26).Range.Delete()

Write-Host "Done"
